# "made changes for month append"
#
# The sheet already has an "Apr-24" Qty/Rate/Value block in columns B:D.
# This appends a matching "May-24" Qty/Rate/Value block in columns E:G
# (same header layout, one data row per existing row 3-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell.
# A plain   $range.Value = "92.05"   lets Excel auto-coerce a numeric-looking
# string into a real number, which is not what we want for the "Rate" column
# (it must stay text, like the existing "Rate" column C). Prefixing with an
# apostrophe forces text, but Excel then stamps the cell with a "quote
# prefix" style (the little green corner-triangle marker) - so we reset the
# cell back to the Normal style right after, leaving a plain text cell with
# no extra formatting, matching the rest of the sheet.
function Set-TextValue($range, [string]$text) {
    $range.Value = '''' + $text
    $range.Style = "Normal"
}

# --- Row 1-2 headers for the new block ---
Set-TextValue $ws.Range("E1") "May-24"
Set-TextValue $ws.Range("F1") ""
Set-TextValue $ws.Range("G1") ""

$ws.Range("E2").Value = "Qty"
$ws.Range("F2").Value = "Rate"
$ws.Range("G2").Value = "Value"

# --- Data rows (Qty / Rate / Value for May-24) ---
$ws.Range("E3").Value = 571087
Set-TextValue $ws.Range("F3") "177.29"
$ws.Range("G3").Value = 101246022

$ws.Range("E4").Value = 142379.9
Set-TextValue $ws.Range("F4") "92.05"
$ws.Range("G4").Value = 13105538.64

$ws.Range("E5").Value = 7416.25
Set-TextValue $ws.Range("F5") "136.27"
$ws.Range("G5").Value = 1010633.4

$ws.Range("E6").Value = 24414.52
Set-TextValue $ws.Range("F6") "190.45"
$ws.Range("G6").Value = 4649643.27

$ws.Range("E7").Value = 563965
Set-TextValue $ws.Range("F7") "167.74"
$ws.Range("G7").Value = 94599326

$ws.Range("E8").Value = 168227.17
Set-TextValue $ws.Range("F8") "195.86"
$ws.Range("G8").Value = 32948492

$ws.Range("E9").Value = 4134
Set-TextValue $ws.Range("F9") "10.00"
$ws.Range("G9").Value = 41340

Set-TextValue $ws.Range("E10") ""
Set-TextValue $ws.Range("F10") ""
$ws.Range("G10").Value = 10.94

Set-TextValue $ws.Range("E11") ""
Set-TextValue $ws.Range("F11") ""
$ws.Range("G11").Value = 4227041

$ws.Range("E12").Value = 8271.5
Set-TextValue $ws.Range("F12") "74.55"
$ws.Range("G12").Value = 616661

$ws.Range("E13").Value = 52769.8
Set-TextValue $ws.Range("F13") "78.99"
$ws.Range("G13").Value = 4168133

Set-TextValue $ws.Range("E14") ""
Set-TextValue $ws.Range("F14") ""
$ws.Range("G14").Value = 125000

Set-TextValue $ws.Range("E15") ""
Set-TextValue $ws.Range("F15") ""
$ws.Range("G15").Value = 1229816

Set-TextValue $ws.Range("E16") ""
Set-TextValue $ws.Range("F16") ""
$ws.Range("G16").Value = 9011835

Set-TextValue $ws.Range("E17") ""
Set-TextValue $ws.Range("F17") ""
$ws.Range("G17").Value = 10366651

Set-TextValue $ws.Range("E18") ""
Set-TextValue $ws.Range("F18") ""
$ws.Range("G18").Value = 2789319.370000002

# --- Column widths for the new columns ---
# Columns F and G now hold data; column H is left pre-formatted with the
# same width (mirroring how column E already had its width set before this
# edit populated it with data), ready for a future month's block.
# ColumnWidth uses Excel's "characters" unit; 14.998697916666666 is the
# value that round-trips to the same OOXML <col width="15.83203125"/> used
# by the existing columns B-E.
$newColWidth = 14.998697916666666
$ws.Columns.Item(6).ColumnWidth = $newColWidth
$ws.Columns.Item(7).ColumnWidth = $newColWidth
$ws.Columns.Item(8).ColumnWidth = $newColWidth

# --- Sheet view flag present on the target file (explicit default value) ---
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch {}
